$d = $word.ActiveDocument
$r = $d.Content
$totalFound = 0

$r = $d.Content
$found = $r.Find.Execute('Appendix 6: SWIFT Quantitative Information Sheet and Consent Form: Participants', $true, $false, $false, $false, $false, $true, 1, $false, 'Bylaag 6: SWIFT Kwantitatiewe Inligtingsblad en Toestemmingsvorm: Deelnemers', 2)
if (-not $found) { Write-Host "NOT FOUND [0]: " 'Appendix 6: SWIFT Quantitative Information Sheet and Consent Form: Participants' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('You’re invited to join a study carried out by researchers from the University of Cape Town in South Africa and the University of Oxford in the United Kingdom.', $true, $false, $false, $false, $false, $true, 1, $false, 'Jy word uitgenooi om aan te sluit by ''n studie wat uitgevoer word deur navorsers van die Universiteit van Kaapstad in Suid-Afrika en die Universiteit van Oxford in die Verenigde Koninkryk.', 2)
if (-not $found) { Write-Host "NOT FOUND [1]: " 'You’re invited to join a study carried out by researchers from the University of Cape Town in South Africa and the University of Oxford in the United Kingdom.' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('We’re doing this study to learn about your experience with chatbots developed by Parenting for Lifelong Health (PLH), Clowns Without Borders South Africa (CWBSA), IDEMS International, and UNICEF South Africa. ', $true, $false, $false, $false, $false, $true, 1, $false, 'Ons doen hierdie studie om meer te leer oor jou ervaring met geselsbots wat ontwikkel is deur Parenting for Lifelong Health (PLH), Clowns Without Borders South Africa (CWBSA), IDEMS International, en UNICEF Suid-Afrika. ', 2)
if (-not $found) { Write-Host "NOT FOUND [2]: " 'We’re doing this study to learn about your experience with chatbots developed by Parenting for Lifelong Health (PLH), Clowns Without Borders South Africa (CWBSA), IDEMS International, and UNICEF South Africa. ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Who can join?', $true, $false, $false, $false, $false, $true, 1, $false, 'Wie kan aansluit?', 2)
if (-not $found) { Write-Host "NOT FOUND [3]: " 'Who can join?' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Do I have to join?', $true, $false, $false, $false, $false, $true, 1, $false, 'Moet ek deelneem?', 2)
if (-not $found) { Write-Host "NOT FOUND [4]: " 'Do I have to join?' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('What happens if I join?', $true, $false, $false, $false, $false, $true, 1, $false, 'Wat gebeur as ek aansluit?', 2)
if (-not $found) { Write-Host "NOT FOUND [5]: " 'What happens if I join?' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('If you decide to join, you will need to read the consent form below and answer “Yes” on WhatsApp to the question, “I have read and understand the information, and I give consent to participate in the study”. ', $true, $false, $false, $false, $false, $true, 1, $false, 'As jy besluit om aan te sluit, moet jy die toestemmingsvorm hieronder lees en "Ja" antwoord op WhatsApp op die vraag, “Ek het die inligting gelees en verstaan, en ek gee toestemming om aan die studie deel te neem.” ', 2)
if (-not $found) { Write-Host "NOT FOUND [6]: " 'If you decide to join, you will need to read the consent form below and answer “Yes” on WhatsApp to the question, “I have read and understand the information, and I give consent to participate in the study”. ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('After completing the first survey, you’ll start to get messages from the ParentText chatbot. ', $true, $false, $false, $false, $false, $true, 1, $false, 'Na die voltooiing van die eerste opname, sal jy begin om boodskappe van die ParentText-geselsbot te ontvang. ', 2)
if (-not $found) { Write-Host "NOT FOUND [7]: " 'After completing the first survey, you’ll start to get messages from the ParentText chatbot. ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('The ParentText chatbot is a 5-day course for parents and caregivers with 5 daily lessons on building a good relationship with your child. ', $true, $false, $false, $false, $false, $true, 1, $false, 'Die ParentText-geselsbot is ''n 5-dag kursus vir ouers en versorgers met 5 daaglikse lesse oor die bou van ''n goeie verhouding met jou kind. ', 2)
if (-not $found) { Write-Host "NOT FOUND [8]: " 'The ParentText chatbot is a 5-day course for parents and caregivers with 5 daily lessons on building a good relationship with your child. ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Do I get anything for joining the study? ', $true, $false, $false, $false, $false, $true, 1, $false, 'Kry ek iets vir deelname aan die studie?', 2)
if (-not $found) { Write-Host "NOT FOUND [9]: " 'Do I get anything for joining the study? ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('What happens to my information if I join?', $true, $false, $false, $false, $false, $true, 1, $false, 'Wat gebeur met my inligting as ek aansluit?', 2)
if (-not $found) { Write-Host "NOT FOUND [10]: " 'What happens to my information if I join?' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('What happens to the research results?', $true, $false, $false, $false, $false, $true, 1, $false, 'Wat gebeur met die navorsingsresultate?', 2)
if (-not $found) { Write-Host "NOT FOUND [11]: " 'What happens to the research results?' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Who are some of the study team members?', $true, $false, $false, $false, $false, $true, 1, $false, 'Wie is sommige van die spanlede van die studie?', 2)
if (-not $found) { Write-Host "NOT FOUND [12]: " 'Who are some of the study team members?' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Is daar enige risiko’s om aan hierdie studie aan te sluit?   ', $true, $false, $false, $false, $false, $true, 1, $false, 'Is daar enige risiko’s om aan hierdie studie deel te neem?   ', 2)
if (-not $found) { Write-Host "NOT FOUND [13]: " 'Is daar enige risiko’s om aan hierdie studie aan te sluit?   ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Ons verwag nie enige risiko’s vir jou as jy aan hierdie studie aansluit nie. As enige vrae jou ongemaklik maak, hoef jy dit nie te antwoord nie. ', $true, $false, $false, $false, $false, $true, 1, $false, 'Ons verwag nie enige risiko’s vir jou as jy aan hierdie studie deelneem nie. As enige vrae jou ongemaklik maak, hoef jy dit nie te antwoord nie. ', 2)
if (-not $found) { Write-Host "NOT FOUND [14]: " 'Ons verwag nie enige risiko’s vir jou as jy aan hierdie studie aansluit nie. As enige vrae jou ongemaklik maak, hoef jy dit nie te antwoord nie. ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Who pays for the study?', $true, $false, $false, $false, $false, $true, 1, $false, 'Wie betaal vir die studie?', 2)
if (-not $found) { Write-Host "NOT FOUND [15]: " 'Who pays for the study?' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('This study is part of the Global Parenting Initiative, which is funded by the LEGO Foundation, Oak Foundation, the World Childhood Foundation, The Human Safety Net, and the UK Research and Innovation Global Challenges Research Fund. ', $true, $false, $false, $false, $false, $true, 1, $false, 'Hierdie studie is deel van die Global Parenting Initiative, gefinansier deur die LEGO Foundation, Oak Foundation, die World Childhood Foundation, The Human Safety Net, en die UK Research and Innovation Global Challenges Research Fund. ', 2)
if (-not $found) { Write-Host "NOT FOUND [16]: " 'This study is part of the Global Parenting Initiative, which is funded by the LEGO Foundation, Oak Foundation, the World Childhood Foundation, The Human Safety Net, and the UK Research and Innovation Global Challenges Research Fund. ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Data protection', $true, $false, $false, $false, $false, $true, 1, $false, 'Databeskerming', 2)
if (-not $found) { Write-Host "NOT FOUND [17]: " 'Data protection' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Who has approved this study?', $true, $false, $false, $false, $false, $true, 1, $false, 'Wie het hierdie studie goedgekeur?', 2)
if (-not $found) { Write-Host "NOT FOUND [18]: " 'Who has approved this study?' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Who do I contact if I have questions or concerns?', $true, $false, $false, $false, $false, $true, 1, $false, 'Wie kan ek kontak as ek vrae of bekommernisse het?', 2)
if (-not $found) { Write-Host "NOT FOUND [19]: " 'Who do I contact if I have questions or concerns?' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('If you have any questions or concerns about your rights as a study participant, you can contact the study team at ', $true, $false, $false, $false, $false, $true, 1, $false, 'As jy enige vrae of bekommernisse het oor jou regte as ''n studie-deelnemer, kan jy die studiespan kontak by ', 2)
if (-not $found) { Write-Host "NOT FOUND [20]: " 'If you have any questions or concerns about your rights as a study participant, you can contact the study team at ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute(' or on WhatsApp at +27 XX XXX XXXX (messages only). ', $true, $false, $false, $false, $false, $true, 1, $false, ' of via WhatsApp by +27 XX XXX XXXX (net boodskappe).', 2)
if (-not $found) { Write-Host "NOT FOUND [21]: " ' or on WhatsApp at +27 XX XXX XXXX (messages only). ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('If you have more questions or concerns about your rights, you can contact this ethics committees: ', $true, $false, $false, $false, $false, $true, 1, $false, 'As jy meer vrae of bekommernisse het oor jou regte, kan jy hierdie etiekkomitees kontak: ', 2)
if (-not $found) { Write-Host "NOT FOUND [22]: " 'If you have more questions or concerns about your rights, you can contact this ethics committees: ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Name', $true, $false, $false, $false, $false, $true, 1, $false, 'Naam', 2)
if (-not $found) { Write-Host "NOT FOUND [23]: " 'Name' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Telephone', $true, $false, $false, $false, $false, $true, 1, $false, 'Telefoon', 2)
if (-not $found) { Write-Host "NOT FOUND [24]: " 'Telephone' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Email', $true, $false, $false, $false, $false, $true, 1, $false, 'E-pos', 2)
if (-not $found) { Write-Host "NOT FOUND [25]: " 'Email' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('University of Cape Town Centre for Social Science Research ', $true, $false, $false, $false, $false, $true, 1, $false, 'Universiteit van Kaapstad Sentrum vir Sosiale Wetenskap Navorsing ', 2)
if (-not $found) { Write-Host "NOT FOUND [26]: " 'University of Cape Town Centre for Social Science Research ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Human Research Ethics Committee', $true, $false, $false, $false, $false, $true, 1, $false, 'Etiekkomitee vir Menslike Navorsing', 2)
if (-not $found) { Write-Host "NOT FOUND [27]: " 'Human Research Ethics Committee' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Informed Consent to Take Part in the Study', $true, $false, $false, $false, $false, $true, 1, $false, 'Ingeligte Toestemming om aan die Studie Deel te Neem', 2)
if (-not $found) { Write-Host "NOT FOUND [28]: " 'Informed Consent to Take Part in the Study' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('Please read these statements carefully: ', $true, $false, $false, $false, $false, $true, 1, $false, 'Lees asseblief hierdie stellings sorgvuldig: ', 2)
if (-not $found) { Write-Host "NOT FOUND [29]: " 'Please read these statements carefully: ' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('I have read the information above and know what I need to do.', $true, $false, $false, $false, $false, $true, 1, $false, 'Ek het die bogenoemde inligting gelees en weet wat ek moet doen.', 2)
if (-not $found) { Write-Host "NOT FOUND [30]: " 'I have read the information above and know what I need to do.' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('I know who can see my information, how it will be kept safe, and what happens to it after the study.', $true, $false, $false, $false, $false, $true, 1, $false, 'Ek weet wie my inligting kan sien, hoe dit veilig gehou sal word, en wat daarmee sal gebeur na die studie.', 2)
if (-not $found) { Write-Host "NOT FOUND [31]: " 'I know who can see my information, how it will be kept safe, and what happens to it after the study.' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('I know I can request access to my data, correct any mistakes, ask to delete it, or for it to be transferred somewhere else.', $true, $false, $false, $false, $false, $true, 1, $false, 'Ek weet ek kan toegang tot my data versoek, enige foute regstel, vra om dit te verwyder, of dit na ''n ander plek oorgedra te word.', 2)
if (-not $found) { Write-Host "NOT FOUND [32]: " 'I know I can request access to my data, correct any mistakes, ask to delete it, or for it to be transferred somewhere else.' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('I know I won’t be named in any papers or reports from this study.', $true, $false, $false, $false, $false, $true, 1, $false, 'Ek weet dat ek nie in enige artikels of verslae van hierdie studie genoem sal word nie.', 2)
if (-not $found) { Write-Host "NOT FOUND [33]: " 'I know I won’t be named in any papers or reports from this study.' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('I know who to contact if I have a problem with the study.', $true, $false, $false, $false, $false, $true, 1, $false, 'Ek weet wie ek kan kontak as ek ''n probleem met die studie het.', 2)
if (-not $found) { Write-Host "NOT FOUND [34]: " 'I know who to contact if I have a problem with the study.' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('You can contact me again if more information is needed from me.', $true, $false, $false, $false, $false, $true, 1, $false, 'Jy kan my weer kontak as daar meer inligting van my benodig word.', 2)
if (-not $found) { Write-Host "NOT FOUND [35]: " 'You can contact me again if more information is needed from me.' } else { $totalFound = $totalFound + 1 }

$r = $d.Content
$found = $r.Find.Execute('You can keep my contact information safe so you can tell me about the results of the study.', $true, $false, $false, $false, $false, $true, 1, $false, 'Jy kan my kontakbesonderhede veilig hou sodat jy my oor die resultate van die studie kan inlig.', 2)
if (-not $found) { Write-Host "NOT FOUND [36]: " 'You can keep my contact information safe so you can tell me about the results of the study.' } else { $totalFound = $totalFound + 1 }

Write-Host "Total replacements applied: $totalFound / 37"